# Database Change in Subject. Sql Code Included.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Apply the "Calculation" cell style to the Religious Studies rows ---
# (rows 7-10, columns B:D) in the first table...
$ws.Range("B7:D10").Style = "Calculation"
# ...and the matching rows (31-34) in the second table.
$ws.Range("B31:D34").Style = "Calculation"

# --- Swap the "Physical Education" / "Work & Life Oriented Education" rows ---
# Row 12 currently holds "Work & Life Oriented Education (155)" and
# row 13 holds "Physical  Education Health  & Sports (147)"; the edit
# swaps their order so row 12 becomes "Physical Education..." and row 13
# becomes "Work & Life...".
$row12 = $ws.Range("B12:D12").Value2
$row13 = $ws.Range("B13:D13").Value2
$ws.Range("B12:D12").Value2 = $row13
$ws.Range("B13:D13").Value2 = $row12

# --- View changes: zoom out to 70% and move the selection to D46 ---
[void]$ws.Select()
[void]$ws.Range("D46").Select()
$excel.ActiveWindow.Zoom = 70
